$wb = $excel.ActiveWorkbook

# --- "Add Line" sheet: append new product rows and tweak row 2 ---
$wsAddLine = $wb.Worksheets.Item("Add Line")

# Update existing row 2 values (Asynchronous Processing -> False, Product -> new id, Unit Price -> 200)
$wsAddLine.Cells.Item(2, 2).Value = $false
$wsAddLine.Cells.Item(2, 7).Value = "a5N0W000001HsIU"
$wsAddLine.Cells.Item(2, 9).Value = 200

# Copy row 2 formatting/values down to rows 3-6, then overwrite the per-row fields
$wsAddLine.Range("A2:K2").Copy($wsAddLine.Range("A3:K3"))
$wsAddLine.Range("A2:K2").Copy($wsAddLine.Range("A4:K4"))
$wsAddLine.Range("A2:K2").Copy($wsAddLine.Range("A5:K5"))
$wsAddLine.Range("A2:K2").Copy($wsAddLine.Range("A6:K6"))

$wsAddLine.Cells.Item(3, 7).Value = "a5N0W000001HpIt"
$wsAddLine.Cells.Item(3, 8).Value = 4
$wsAddLine.Cells.Item(3, 9).Value = 100
$wsAddLine.Cells.Item(3, 10).Value = 2

$wsAddLine.Cells.Item(4, 7).Value = "a5N6T0000011eQy"
$wsAddLine.Cells.Item(4, 8).Value = 2
$wsAddLine.Cells.Item(4, 9).Value = 50
$wsAddLine.Cells.Item(4, 10).Value = 3

$wsAddLine.Cells.Item(5, 7).Value = "a5N0W000001I396"
$wsAddLine.Cells.Item(5, 8).Value = 2
$wsAddLine.Cells.Item(5, 9).Value = 80
$wsAddLine.Cells.Item(5, 10).Value = 4

$wsAddLine.Cells.Item(6, 7).Value = "a5N0W000001I39B"
$wsAddLine.Cells.Item(6, 8).Value = 3
$wsAddLine.Cells.Item(6, 9).Value = 20
$wsAddLine.Cells.Item(6, 10).Value = 5

# --- "Add Header" sheet: move selection only ---
$wsAddHeader = $wb.Worksheets.Item("Add Header")
$wsAddHeader.Range("E14").Select()

# --- "Add Line" becomes the active/selected sheet & cell ---
$wsAddLine.Activate()
$wsAddLine.Range("C16").Select()
